$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "244.17"
Set-TextValue $ws.Range("G2") "15"

Set-TextValue $ws.Range("D3") "23.90"
Set-TextValue $ws.Range("G3") "15"

Set-TextValue $ws.Range("D4") "5.281"
Set-TextValue $ws.Range("G4") "15"

Set-TextValue $ws.Range("D5") "0.05882"
Set-TextValue $ws.Range("G5") "15"

Set-TextValue $ws.Range("D6") "6.475"
Set-TextValue $ws.Range("G6") "15"

Set-TextValue $ws.Range("D7") "3.335"
Set-TextValue $ws.Range("G7") "15"

Set-TextValue $ws.Range("D8") "0.8145"
Set-TextValue $ws.Range("G8") "15"

Set-TextValue $ws.Range("D9") "0.8949"
Set-TextValue $ws.Range("G9") "15"

Set-TextValue $ws.Range("D10") "0.1387"
Set-TextValue $ws.Range("G10") "15"

Set-TextValue $ws.Range("D11") "0.07238"
Set-TextValue $ws.Range("G11") "15"

Set-TextValue $ws.Range("D12") "0.03080"
Set-TextValue $ws.Range("G12") "15"

Set-TextValue $ws.Range("D13") "0.03032"
Set-TextValue $ws.Range("G13") "15"

Set-TextValue $ws.Range("D14") "0.09348"
Set-TextValue $ws.Range("G14") "15"

Set-TextValue $ws.Range("D15") "3.848"
Set-TextValue $ws.Range("G15") "15"

Set-TextValue $ws.Range("D16") "0.001550"
Set-TextValue $ws.Range("G16") "15"

Set-TextValue $ws.Range("D17") "0.04698"
Set-TextValue $ws.Range("G17") "15"

Set-TextValue $ws.Range("D18") "0.0006004"
Set-TextValue $ws.Range("E18") "17OneONE"
Set-TextValue $ws.Range("G18") "15"

Set-TextValue $ws.Range("D19") "0.006241"
Set-TextValue $ws.Range("G19") "15"

Set-TextValue $ws.Range("D20") "0.001262"
Set-TextValue $ws.Range("G20") "15"

Set-TextValue $ws.Range("D21") "0.004613"
Set-TextValue $ws.Range("G21") "15"

Set-TextValue $ws.Range("D22") "0.00008693"
Set-TextValue $ws.Range("G22") "15"

Set-TextValue $ws.Range("D23") "3.562"
Set-TextValue $ws.Range("G23") "15"

Set-TextValue $ws.Range("D24") "2.180"
Set-TextValue $ws.Range("G24") "15"

Set-TextValue $ws.Range("D25") "0.3205"
Set-TextValue $ws.Range("G25") "15"

Set-TextValue $ws.Range("G26") "15"

Set-TextValue $ws.Range("G27") "15"

Set-TextValue $ws.Range("D28") "0.0002337"
Set-TextValue $ws.Range("G28") "15"

Set-TextValue $ws.Range("G29") "15"

Set-TextValue $ws.Range("G30") "15"

Set-TextValue $ws.Range("G31") "15"

Set-TextValue $ws.Range("G32") "15"

Set-TextValue $ws.Range("G33") "15"

Set-TextValue $ws.Range("G34") "15"

Set-TextValue $ws.Range("G35") "15"

Set-TextValue $ws.Range("G36") "15"

Set-TextValue $ws.Range("G37") "15"

Set-TextValue $ws.Range("G38") "15"

Set-TextValue $ws.Range("G39") "15"

Set-TextValue $ws.Range("D40") "0.03799"
Set-TextValue $ws.Range("G40") "15"

Set-TextValue $ws.Range("D41") "0.006345"
Set-TextValue $ws.Range("G41") "15"

Set-TextValue $ws.Range("D42") "0.1058"
Set-TextValue $ws.Range("G42") "15"

Set-TextValue $ws.Range("D43") "0.002698"
Set-TextValue $ws.Range("G43") "15"

Set-TextValue $ws.Range("D44") "0.007092"
Set-TextValue $ws.Range("G44") "15"

Set-TextValue $ws.Range("D45") "0.00005399"
Set-TextValue $ws.Range("G45") "15"

Set-TextValue $ws.Range("G46") "15"

Set-TextValue $ws.Range("D47") "0.5395"
Set-TextValue $ws.Range("E47") "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue $ws.Range("G47") "15"

Set-TextValue $ws.Range("D48") "0.02104"
Set-TextValue $ws.Range("G48") "15"

Set-TextValue $ws.Range("D49") "0.00002098"
Set-TextValue $ws.Range("G49") "15"

Set-TextValue $ws.Range("D50") "0.0001998"
Set-TextValue $ws.Range("G50") "15"

Set-TextValue $ws.Range("G51") "15"
